$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 (Mississippi) - fill in the previously-empty/error values with results
$ws.Range("B42").Value = 44022
$ws.Range("C42").Value = 35419
$ws.Range("D42").Value = 1230
$ws.Range("E42").Value = 16811
$ws.Range("F42").Value = 616
$ws.Range("G42").Value = 47.46
$ws.Range("H42").Value = 50.08

# Copy the date number format from a neighboring "Date Published" cell (B41)
# so B42 matches style s="2" (numFmt 165 "YYYY-MM-DD")
$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial(-4122) | Out-Null
$ws.Range("B42").Value = 44022
$excel.CutCopyMode = $false

$ws.Range("I42").Value = $true
$ws.Range("J42").Value = $true

$ws.Range("O42").Value = "Success!"
